# Comments updated in excel sheet
# - Add a new "Steps" explanatory note in cell C16 of the "Steps" sheet,
#   highlighted with a bold font on a yellow fill (new shared string +
#   new fill/cellXf in styles.xml).
# - Move the active selection on the "Steps" sheet to C17.
# - Set the Steps sheet's page setup (paper size / orientation) so a
#   pageSetup element is written for that sheet, matching the other sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Steps")

# New note cell, bold text on a yellow highlight fill.
$c16 = $ws.Range("C16")
$c16.Value = "What to modify in code? -  What parameter is on what line in the code as what variable?"
$c16.Font.Bold = $true
$c16.Interior.Color = 65535

# Page setup for the Steps sheet (adds <pageSetup .../> like the other sheets).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the selection/active cell to C17 as in the authored workbook.
[void]$ws.Range("C17").Select()
